# Renames the diff-report's "old/new" header suffixes to the concrete
# format-version names they actually represent (FV2310 / FV2404), then
# turns the header row into a proper Excel Table and freezes it so it
# stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rewrite the 21 header cells (row 1, columns A:U) ------------------
$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn A1:U88 into a real Table so the headers gain filter buttons --
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U88"), 0, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row so it stays put while scrolling -------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
